$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "pruebaTecnica"

# The two previously-styled cells (C2 quote-prefixed, C3 underline-font)
# must come back to the plain/default style before we repopulate them.
$ws.Cells.Item(2,3).Style = "Normal"
$ws.Cells.Item(3,3).Style = "Normal"

# ---------------------------------------------------------------------------
# Row 1 (headers) - all header cells pick up the existing bold/centered
# row style (s="1") automatically when written.
# ---------------------------------------------------------------------------
$ws.Cells.Item(1,1).Value  = "TAG"
$ws.Cells.Item(1,2).Value  = "keyUserName"
$ws.Cells.Item(1,3).Value  = "keyPassword"
$ws.Cells.Item(1,4).Value  = "keyBook1"
$ws.Cells.Item(1,5).Value  = "keyBook2"
$ws.Cells.Item(1,6).Value  = "keyFirstName"
$ws.Cells.Item(1,7).Value  = "keyLastName"
$ws.Cells.Item(1,8).Value  = "keyEmail"
$ws.Cells.Item(1,9).Value  = "keyGender"
$ws.Cells.Item(1,10).Value = "keyMobile"
$ws.Cells.Item(1,11).Value = "keyBirthDate"
$ws.Cells.Item(1,12).Value = "keySubjects"
$ws.Cells.Item(1,13).Value = "keyHobbies"
$ws.Cells.Item(1,14).Value = "keyCurrentAddress"
$ws.Cells.Item(1,15).Value = "keyState"
$ws.Cells.Item(1,16).Value = "keyCity"

# ---------------------------------------------------------------------------
# Row 2
# ---------------------------------------------------------------------------
$ws.Cells.Item(2,1).Value  = " @PruebaTecnicaPunto1"
$ws.Cells.Item(2,2).Value  = "test01"
$ws.Cells.Item(2,3).Value  = "Test2023*"
$ws.Cells.Item(2,4).Value  = "Programming JavaScript"
$ws.Cells.Item(2,5).Value  = "Understanding ECMAScript 6"
$ws.Cells.Item(2,6).Value  = "-"
$ws.Cells.Item(2,7).Value  = "-"
$ws.Cells.Item(2,8).Value  = "-"
$ws.Cells.Item(2,9).Value  = "-"
$ws.Cells.Item(2,10).Value = "-"

# K2:P2 use a date number format (167 -> "d mmm yyyy") even though they just
# hold a placeholder dash.
$ws.Range("K2:P2").NumberFormat = "d\ mmm\ yyyy"
$ws.Cells.Item(2,11).Value = "-"
$ws.Cells.Item(2,12).Value = "-"
$ws.Cells.Item(2,13).Value = "-"
$ws.Cells.Item(2,14).Value = "-"
$ws.Cells.Item(2,15).Value = "-"
$ws.Cells.Item(2,16).Value = "-"

# ---------------------------------------------------------------------------
# Row 3
# ---------------------------------------------------------------------------
$ws.Cells.Item(3,1).Value = " @PruebaTecnicaPunto2"
$ws.Cells.Item(3,2).Value = "test02"
$ws.Cells.Item(3,3).Value = "Test2023*"
$ws.Cells.Item(3,4).Value = "-"
$ws.Cells.Item(3,5).Value = "-"
$ws.Cells.Item(3,6).Value = "Prueba01"
$ws.Cells.Item(3,7).Value = "Prueba01"
$ws.Cells.Item(3,8).Value = "prueba@yopmail.com"
$ws.Cells.Item(3,9).Value = "Female"

# J3 keeps a text quote-prefix (mobile number stored as text).
$ws.Cells.Item(3,10).Value = "'3004441234"

# K3 is quote-prefixed text too, but also carries the date number format.
$ws.Cells.Item(3,11).NumberFormat = "d\ mmm\ yyyy"
$ws.Cells.Item(3,11).Value = "'27 Sep 2000"

$ws.Cells.Item(3,12).Value = "Computer"
$ws.Cells.Item(3,13).Value = "Music"
$ws.Cells.Item(3,14).Value = "Autopista Norte al oriente"
$ws.Cells.Item(3,15).Value = "Haryana"
$ws.Cells.Item(3,16).Value = "karnal"

# ---------------------------------------------------------------------------
# Row 4 - single empty, underlined placeholder cell at L4 (re-using the
# underline font already present in the workbook, same as old C3).
# ---------------------------------------------------------------------------
$ws.Cells.Item(4,12).Font.Underline = 2

# ---------------------------------------------------------------------------
# Column widths (best effort - engine quantizes ColumnWidth to 1/6 character
# increments, so these land as close as possible to the authored widths).
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 12.9
$ws.Range("F1:G1").ColumnWidth = 12.9
$ws.Columns.Item(11).ColumnWidth = 10.9

# ---------------------------------------------------------------------------
# View state: freeze-pane stays the same; move the active selection to L4,
# matching the saved selection in the edited file.
# ---------------------------------------------------------------------------
$ws.Range("L4").Select()
